$d = $word.ActiveDocument

# 1. Merge the split "Tut" + "o" + "rials" runs back into a single
#    "AVR I/O Register Configuration - Tutorials" run inside the hyperlink.
$d.Content.Find.Execute("AVR I/O Register Configuration - Tutorials", $false, $false, $false, $false, $false, $true, 1, $false, "AVR I/O Register Configuration - Tutorials", 2)
